# Update "想去人数" (F column) values on both the "展览" and "全部类型"
# worksheets, which contain duplicated data rows.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new F-column value
$updates = @{
    2  = 631
    5  = 13158
    6  = 75
    13 = 14377
    25 = 5430
    26 = 939
    28 = 318
    30 = 54
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
